$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "           0       0.96      0.87      0.91      8889"
$ws.Range("A4").Value = "           1       0.88      0.97      0.92      8897"
$ws.Range("A6").Value = "    accuracy                           0.92     17786"
$ws.Range("A7").Value = "   macro avg       0.92      0.92      0.92     17786"
$ws.Range("A8").Value = "weighted avg       0.92      0.92      0.92     17786"

# A13/E13 start with "[" / end with "]" so Excel will not coerce them to numbers,
# but B13/C13/D13 look like plain (space-prefixed) numbers and must be forced to
# stay text so the leading space and string type are preserved.
$ws.Range("A13").Value = "[0.8054202198982239"

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = " 0.864345133304596"

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = " 0.894369900226593"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = " 0.9108066558837891"

$ws.Range("E13").Value = " 0.9222955107688904]"
